# Updated symbol list on Tue Jan 24 11:49:43 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) columns for the crypto symbol table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = <new price text>; E = <new volume% text> } (omitted key = unchanged)
$updates = @{
    2 = @{ D="315.23"; E="2.91%" }
    3 = @{ D="35.33"; E="-2.19%" }
    4 = @{ D="5.132"; E="0.18%" }
    6 = @{ D="2.128"; E="-0.14%" }
    7 = @{ D="4.155"; E="0.84%" }
    8 = @{ E="-0.08%" }
    9 = @{ D="0.9292"; E="0.99%" }
    10 = @{ D="0.1013"; E="4.47%" }
    11 = @{ D="0.1869"; E="0.56%" }
    12 = @{ D="0.09105"; E="4.21%" }
    13 = @{ D="0.03608"; E="1.05%" }
    14 = @{ D="0.09906"; E="-0.10%" }
    15 = @{ D="0.001443"; E="0.57%" }
    16 = @{ D="0.005689"; E="0.12%" }
    17 = @{ E="-0.11%" }
    18 = @{ D="2.890"; E="9.82%" }
    19 = @{ D="0.3411"; E="0.58%" }
    20 = @{ E="-0.50%" }
    21 = @{ D="5.106"; E="-1.19%" }
    22 = @{ D="0.2221"; E="9.92%" }
    23 = @{ D="0.04556"; E="-0.25%" }
    24 = @{ D="0.001246"; E="0.93%" }
    25 = @{ D="0.004700"; E="-6.69%" }
    26 = @{ E="-21.93%" }
    27 = @{ D="0.0004505"; E="-5.14%" }
    39 = @{ D="0.01960"; E="5.45%" }
    40 = @{ D="0.04851"; E="1.91%" }
    41 = @{ D="0.007726"; E="1.22%" }
    42 = @{ E="-0.40%" }
    43 = @{ D="0.007846"; E="1.24%" }
    44 = @{ D="0.002152"; E="-2.66%" }
    45 = @{ E="6.97%" }
    46 = @{ D="0.00006649"; E="5.06%" }
    47 = @{ E="0.11%" }
    48 = @{ D="39.45"; E="-18.14%" }
    49 = @{ D="0.001702"; E="-14.90%" }
    50 = @{ E="0.11%" }
    51 = @{ E="0.11%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    foreach ($col in $vals.Keys) {
        $cell = $ws.Range("$col$row")
        # Force text storage so numeric-looking strings (e.g. "315.23", "2.91%")
        # are not silently reinterpreted as Number/Percentage values.
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col]
        $cell.Style = "Normal"
    }
}
